{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies a series of unique text search-and-replace operations that\n// together reproduce the target diff (job title, and several small\n// wording tweaks in the body paragraphs of the cover letter).\n\nconst body = context.document.body;\n\n// Each entry: [uniqueSearchText, replacementText]\nconst replacements = [\n  // \"...at Lockheed Martin.\" -> \"...at Nvidia.\"\n  [\"Lockheed Martin\", \"Nvidia\"],\n\n  // \"My education, my ...\" -> \"Given my education, my ...\"\n  [\"My education, my \", \"Given my education, my \"],\n\n  // \"...and my passion for technology make an ideal candidate...\"\n  // -> \"...and my passion for technology I believe I am an ideal candidate...\"\n  [\", and my passion for technology make \", \", and my passion for technology I believe I am \"],\n\n  // \"I come across fascinating and revolutionary technology\"\n  // -> \"I interact with fascinating and revolutionary technology\"\n  [\n    \"As an engineering student working on various projects, I come across fascinating and \",\n    \"As an engineering student working on various projects, I interact with fascinating and \",\n  ],\n\n  // \"interested in connecting clients and their needs.\"\n  // -> \"interested in connecting these clients to their needs.\"\n  [\"connecting clients and their needs\", \"connecting these clients to their needs\"],\n\n  // \"...a TN visa as a Canadian professional and citizen. \"\n  // -> \"...a TN visa sponsorship as a Canadian professional. \"\n  [\" as a Canadian professional and citizen. \", \" sponsorship as a Canadian professional. \"],\n];\n\nfor (const [searchText, replacementText] of replacements) {\n  const found = body.search(searchText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${searchText}`);\n  }\n\n  found.items[0].insertText(replacementText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies a series of unique Find/Replace operations that together\n# reproduce the target diff (job title, and several small wording\n# tweaks in the body paragraphs of the cover letter).\n\n$doc = $word.ActiveDocument\n\n$wdReplaceAll = \"wdReplaceAll\"\n$wdFindContinue = 1\n\nfunction Replace-UniqueText($FindText, $ReplaceText) {\n    $range = $doc.Content\n    $found = $range.Find.Execute(\n        $FindText,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        $wdFindContinue,\n        $false,\n        $ReplaceText,\n        $wdReplaceAll\n    )\n\n    if (-not $found) {\n        throw \"Text not found: $FindText\"\n    }\n}\n\n# \"...at Lockheed Martin.\" -> \"...at Nvidia.\"\nReplace-UniqueText \"Lockheed Martin\" \"Nvidia\"\n\n# \"My education, my ...\" -> \"Given my education, my ...\"\nReplace-UniqueText \"My education, my \" \"Given my education, my \"\n\n# \"...and my passion for technology make an ideal candidate...\"\n# -> \"...and my passion for technology I believe I am an ideal candidate...\"\nReplace-UniqueText \", and my passion for technology make \" \", and my passion for technology I believe I am \"\n\n# \"I come across fascinating and revolutionary technology\"\n# -> \"I interact with fascinating and revolutionary technology\"\nReplace-UniqueText \"As an engineering student working on various projects, I come across fascinating and \" \"As an engineering student working on various projects, I interact with fascinating and \"\n\n# \"interested in connecting clients and their needs.\"\n# -> \"interested in connecting these clients to their needs.\"\nReplace-UniqueText \"connecting clients and their needs\" \"connecting these clients to their needs\"\n\n# \"...a TN visa as a Canadian professional and citizen. \"\n# -> \"...a TN visa sponsorship as a Canadian professional. \"\nReplace-UniqueText \" as a Canadian professional and citizen. \" \" sponsorship as a Canadian professional. \"\n\nWrite-Output \"done\"\n"}
